# Renames the Pearson / BTEC logo picture shapes that live in the
# document's headers and footers, swapping their "friendly name"
# (image2.png <-> image1.png, image1.jpg <-> image2.jpg) while leaving
# everything else (the embedded picture data, alt text/description,
# size, position, etc.) untouched.
#
# wdHeaderFooterIndex mapping used below:
#   Headers/Footers.Item(1) -> wdHeaderFooterPrimary   (the "default" header/footer)
#   Headers/Footers.Item(2) -> wdHeaderFooterFirstPage (the "first page" header/footer)

$d = $word.ActiveDocument

function Rename-InlineLogo($range, [string]$newName) {
    $shape = $range.InlineShapes.Item(1)
    $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

foreach ($sec in $d.Sections) {

    # Footer, first page (footer1.xml) - Pearson logo: image2.png -> image1.png
    $footerFirst = $sec.Footers.Item(2)
    if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -gt 0) {
        Rename-InlineLogo $footerFirst.Range "image1.png"
    }

    # Footer, default/primary (footer2.xml) - Pearson logo: image2.png -> image1.png
    $footerDefault = $sec.Footers.Item(1)
    if ($footerDefault.Exists -and $footerDefault.Range.InlineShapes.Count -gt 0) {
        Rename-InlineLogo $footerDefault.Range "image1.png"
    }

    # Header, first page (header1.xml) - BTEC logo: image1.jpg -> image2.jpg
    $headerFirst = $sec.Headers.Item(2)
    if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -gt 0) {
        Rename-InlineLogo $headerFirst.Range "image2.jpg"
    }
}
